# The "Gaz" sheet was missing the weekend dates 2025-06-21 (Saturday) and
# 2025-06-22 (Sunday). Insert them before the existing 2025-06-23 row,
# shifting the rows that follow (2025-06-23, 2025-06-25) down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# Remember the rows that need to shift down before we overwrite them.
$oldA7 = $ws.Range("A7").Value()
$oldB7 = $ws.Range("B7").Value()
$oldA8 = $ws.Range("A8").Value()
$oldB8 = $ws.Range("B8").Value()

# Old row 8 (2025-06-25 | 34.75) -> row 10.
$ws.Range("A10").Value = "'" + $oldA8
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = $oldB8

# Old row 7 (2025-06-23 | 40.9) -> row 9.
$ws.Range("A9").Value = "'" + $oldA7
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = $oldB7

# New row 7: 2025-06-21 | 40.275 (Saturday, price carried from source data).
$ws.Range("A7").Value = "'2025-06-21"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 40.275

# New row 8: 2025-06-22 | 40.275 (Sunday, price carried from source data).
$ws.Range("A8").Value = "'2025-06-22"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 40.275
